# Update gh-pages to output generated at 456a3b4
# Applies the refreshed "想去人数" (interest count) figures - and a couple of
# related detail/date/cover refreshes - across all four sheets of the
# 广州-漫展信息 workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 27
$ws.Range("F3").Value = 2764
$ws.Range("F4").Value = 1099
$ws.Range("F5").Value = 20140
$ws.Range("F6").Value = 84
$ws.Range("F7").Value = 2368
$ws.Range("F8").Value = 764
$ws.Range("F10").Value = 463
$ws.Range("F11").Value = 709
$ws.Range("F12").Value = 256
$ws.Range("F15").Value = 385
$ws.Range("F16").Value = 90
$ws.Range("F17").Value = 478
$ws.Range("F19").Value = 223
$ws.Range("F21").Value = 21
$ws.Range("F22").Value = 107

# ---------------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 199
$ws.Range("F6").Value = 122
$ws.Range("F14").Value = 93
$ws.Range("F16").Value = 106
$ws.Range("F20").Value = 1
$ws.Range("F22").Value = 36

# ---------------------------------------------------------------------
# Sheet 3: 本地生活 (Local Life)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 660
$ws.Range("E4").Value = "2024.08.30 00:00-10.31 23:59"
$ws.Range("F4").Value = 608
$ws.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202409/khbZM0d01727406753078.jpeg"
$ws.Range("F5").Value = 1130

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (All Types) - union of the sheets above, sorted by date
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 660
$ws.Range("E4").Value = "2024.08.30 00:00-10.31 23:59"
$ws.Range("F4").Value = 608
$ws.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202409/khbZM0d01727406753078.jpeg"
$ws.Range("F5").Value = 199
$ws.Range("F6").Value = 27
$ws.Range("F8").Value = 1131
$ws.Range("F9").Value = 2764
$ws.Range("F10").Value = 1099
$ws.Range("F11").Value = 20140
$ws.Range("F14").Value = 84
$ws.Range("F15").Value = 122
$ws.Range("F17").Value = 2368
$ws.Range("F18").Value = 764
$ws.Range("F21").Value = 463
$ws.Range("F22").Value = 709
$ws.Range("F23").Value = 256
$ws.Range("F29").Value = 385
$ws.Range("F30").Value = 90
$ws.Range("F33").Value = 478
$ws.Range("F34").Value = 93
$ws.Range("F37").Value = 223
$ws.Range("F38").Value = 106
$ws.Range("F39").Value = 106
$ws.Range("F44").Value = 21
$ws.Range("F45").Value = 1
$ws.Range("F47").Value = 36
$ws.Range("F50").Value = 107
